# TestSuit.xlsx update — "config file added and util package added"
#
# Net change vs before.xlsx:
#   - Execution Flag (column C) flips from "NO" to "YES" for rows 4-12
#     (these scenarios are now covered now that the config file / util
#     package landed).
#   - A brand new scenario row (row 19) is appended:
#       A19 = "checking checkout page"
#       B19 = "\u201c\u201d"
#       C19 = "NO"
#       D19 = "login"
#       E19 = "checkout"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flip Execution Flag NO -> YES for rows 4 through 12.
$ws.Range("C4").Value = "YES"
$ws.Range("C5").Value = "YES"
$ws.Range("C6").Value = "YES"
$ws.Range("C7").Value = "YES"
$ws.Range("C8").Value = "YES"
$ws.Range("C9").Value = "YES"
$ws.Range("C10").Value = "YES"
$ws.Range("C11").Value = "YES"
$ws.Range("C12").Value = "YES"

# Append the new test scenario as row 19.
$openQuote = [char]0x201C
$closeQuote = [char]0x201D
$ws.Range("A19").Value = "checking checkout page"
$ws.Range("B19").Value = "$openQuote$closeQuote"
$ws.Range("C19").Value = "NO"
$ws.Range("D19").Value = "login"
$ws.Range("E19").Value = "checkout"

# Selection ends up parked past the new data, matching the authored file.
$ws.Range("E21").Select()
